$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44263
$ws.Range("M2").Value = 250
$ws.Range("N2").Value = 21000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 21500
$ws.Range("Q2").Value = '$/caja 18 kilos'
$ws.Range("S2").Value = 1194
$ws.Range("T2").Value = 18

$ws.Range("D3").Value = 44489
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 26000
$ws.Range("O3").Value = 27000
$ws.Range("P3").Value = 26500
$ws.Range("Q3").Value = '$/bandeja 18 kilos'
$ws.Range("S3").Value = 1472
$ws.Range("T3").Value = 18

$ws.Range("D4").Value = 44291
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("Q4").Value = '$/bandeja 18 kilos'
$ws.Range("S4").Value = 972
$ws.Range("T4").Value = 18

$ws.Range("D5").Value = 44418
$ws.Range("M5").Value = 240
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 11000
$ws.Range("P5").Value = 10500
$ws.Range("Q5").Value = '$/bandeja 10 kilos'
$ws.Range("S5").Value = 1050
$ws.Range("T5").Value = 10

$ws.Range("D6").Value = 44323
$ws.Range("M6").Value = 270
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("Q6").Value = '$/bandeja 18 kilos'
$ws.Range("S6").Value = 1194
$ws.Range("T6").Value = 18

$ws.Range("D7").Value = 44487
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 14500
$ws.Range("Q7").Value = '$/bandeja 10 kilos'
$ws.Range("S7").Value = 1450
$ws.Range("T7").Value = 10

$ws.Range("D8").Value = 44307
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 19000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 19500
$ws.Range("Q8").Value = '$/bandeja 18 kilos'
$ws.Range("S8").Value = 1083
$ws.Range("T8").Value = 18
